$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4266
$ws1.Range("F3").Value = 2423
$ws1.Range("C4").Value = "南宁·恋与深空only（取消）"
$ws1.Range("F4").Value = 480
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F9").Value = 123
$ws1.Range("F10").Value = 132
$ws1.Range("F11").Value = 151
$ws1.Range("F13").Value = 291
$ws1.Range("F14").Value = 3283
$ws1.Range("F15").Value = 222

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4266
$ws4.Range("F3").Value = 2423
$ws4.Range("C4").Value = "南宁·恋与深空only（取消）"
$ws4.Range("F4").Value = 480
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F6").Value = 1
$ws4.Range("F11").Value = 123
$ws4.Range("F12").Value = 132
$ws4.Range("F13").Value = 151
$ws4.Range("F17").Value = 291
$ws4.Range("F18").Value = 3283
$ws4.Range("F19").Value = 222
